$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes: new column D (English label mirrored from/duplicated onto the
# translated C column), localized Malay values moved into C for several lookups,
# and six new "penalty_price" rows (57-62) with columns A-F.
$ws.Range("D1").Value = 'Johor'
$ws.Range("D2").Value = 'Kedah'
$ws.Range("D3").Value = 'Kelantan'
$ws.Range("D4").Value = 'Melaka'
$ws.Range("D5").Value = 'Negeri Sembilan'
$ws.Range("D6").Value = 'Pahang'
$ws.Range("D7").Value = 'Penang'
$ws.Range("D8").Value = 'Perak'
$ws.Range("D9").Value = 'Perlis'
$ws.Range("D10").Value = 'Selangor'
$ws.Range("D11").Value = 'Terengganu'
$ws.Range("D12").Value = 'Sabah'
$ws.Range("D13").Value = 'Sarawak'
$ws.Range("D14").Value = 'Wilayah Persekutuan Kuala Lumpur'
$ws.Range("D15").Value = 'Wilayah Persekutuan Labuan'
$ws.Range("D16").Value = 'Wilayah Persekutuan Putrajaya'
$ws.Range("D17").Value = 'Major'
$ws.Range("D18").Value = 'TBPS'
$ws.Range("D19").Value = 'DISC'
$ws.Range("D20").Value = 'Isu Dalaman'
$ws.Range("D21").Value = 'Isu Luaran'
$ws.Range("D22").Value = 'Aktif'
$ws.Range("D23").Value = 'Tidak Aktif'
$ws.Range("D24").Value = 'Kritikal'
$ws.Range("D25").Value = 'Tidak Kritikal'
$ws.Range("D26").Value = 'E-mel'
$ws.Range("D27").Value = 'Sembang Langsung'
$ws.Range("D28").Value = 'PTJ'
$ws.Range("C29").Value = 'branch'
$ws.Range("D29").Value = 'Cawangan'
$ws.Range("D30").Value = 'Isnin'
$ws.Range("D31").Value = 'Selasa'
$ws.Range("D32").Value = 'Rabu'
$ws.Range("D33").Value = 'Khamis'
$ws.Range("D34").Value = 'Jumaat'
$ws.Range("D35").Value = 'Sabtu'
$ws.Range("D36").Value = 'Ahad'
$ws.Range("D37").Value = 'Normal'
$ws.Range("C38").Value = 'Half Day'
$ws.Range("D38").Value = 'Separuh Hari'
$ws.Range("C39").Value = 'Weekend'
$ws.Range("D39").Value = 'Hujung Minggu'
$ws.Range("C40").Value = 'To: Complainant, CC: Relevant Group, BCC: Not Applicable'
$ws.Range("D40").Value = 'To: Pengadu, CC: Kumpulan yang berkaitan BCC: Tidak Berkenaan'
$ws.Range("C41").Value = 'To: Relevant Technician, CC: Relevant Technician Group, BCC: Complainant'
$ws.Range("D41").Value = 'To: Juruteknik yang berkaitan, CC: Kumpulan Juruteknik yang berkaitan BC: Pengadu'
$ws.Range("C42").Value = 'To: Relevant Technician, CC: Relevant Technician Group, BCC: Not Applicable'
$ws.Range("D42").Value = 'To: Juruteknik yang berkaitan, CC: Kumpulan Juruteknik yang berkaitan BCC: Tidak Berkenaan'
$ws.Range("C43").Value = 'Minutes'
$ws.Range("D43").Value = 'Minit'
$ws.Range("C44").Value = 'Hours'
$ws.Range("D44").Value = 'Jam'
$ws.Range("C45").Value = 'Days'
$ws.Range("D45").Value = 'Hari'
$ws.Range("C46").Value = 'Not Important'
$ws.Range("D46").Value = 'Tidak Penting '
$ws.Range("C47").Value = 'Critical'
$ws.Range("D47").Value = 'Kritikal'
$ws.Range("C48").Value = 'Important'
$ws.Range("D48").Value = 'Penting'
$ws.Range("C49").Value = 'Moderate'
$ws.Range("D49").Value = 'Sederhana'
$ws.Range("C50").Value = 'Low'
$ws.Range("D50").Value = 'Rendah'
$ws.Range("C51").Value = 'Missing'
$ws.Range("D51").Value = 'Hilang'
$ws.Range("C52").Value = 'Damaged / Broken'
$ws.Range("D52").Value = 'Rosak'
$ws.Range("C53").Value = 'Phone'
$ws.Range("D53").Value = 'Telefon'
$ws.Range("C54").Value = 'Email'
$ws.Range("D54").Value = 'Emel'
$ws.Range("D55").Value = 'Chatbot'
$ws.Range("D56").Value = 'Live Chat'
$ws.Range("A57").Value = 'penalty_price'
$ws.Range("B57").Value = 1
$ws.Range("C57").Value = 'RM 10000'
$ws.Range("D57").Value = 'RM 10000'
$ws.Range("E57").Value = 'severity'
$ws.Range("F57").Value = 'Not Important'
$ws.Range("A58").Value = 'penalty_price'
$ws.Range("B58").Value = 2
$ws.Range("C58").Value = 'RM 5000'
$ws.Range("D58").Value = 'RM 5000'
$ws.Range("E58").Value = 'severity'
$ws.Range("F58").Value = 'Critical'
$ws.Range("A59").Value = 'penalty_price'
$ws.Range("B59").Value = 3
$ws.Range("C59").Value = 'RM 3000'
$ws.Range("D59").Value = 'RM 3000'
$ws.Range("E59").Value = 'severity'
$ws.Range("F59").Value = 'Important'
$ws.Range("A60").Value = 'penalty_price'
$ws.Range("B60").Value = 4
$ws.Range("C60").Value = 'RM 1000'
$ws.Range("D60").Value = 'RM 1000'
$ws.Range("E60").Value = 'severity'
$ws.Range("F60").Value = 'Not Important'
$ws.Range("A61").Value = 'penalty_price'
$ws.Range("B61").Value = 5
$ws.Range("C61").Value = 'RM 500'
$ws.Range("D61").Value = 'RM 500'
$ws.Range("E61").Value = 'severity'
$ws.Range("F61").Value = 'Critical'
$ws.Range("A62").Value = 'penalty_price'
$ws.Range("B62").Value = 6
$ws.Range("C62").Value = 'RM 200'
$ws.Range("D62").Value = 'RM 200'
$ws.Range("E62").Value = 'severity'
$ws.Range("F62").Value = 'Important'

# Column D width (matches the new 92.85-ish-style custom column for the English label)
$ws.Columns.Item(4).ColumnWidth = 42.2

# C56 ("Live Chat") gets wrap text + vertical-centered alignment
$c56 = $ws.Range("C56")
$c56.WrapText = $true
$c56.VerticalAlignment = -4108

# Restore a representative selection/scroll position similar to the authored state
$ws.Range("H52").Select()
